$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(76, 8).Value = 19610758
$ws.Cells.Item(76, 9).Value = 3074
$ws.Cells.Item(76, 10).Value = 47621736
$ws.Cells.Item(76, 11).Value = 3074
$ws.Cells.Item(76, 12).Value = 47621736
$ws.Cells.Item(76, 13).Value = -2759
$ws.Cells.Item(76, 14).Value = -47622366

$ws.Cells.Item(79, 8).Value = 19610758
$ws.Cells.Item(79, 9).Value = 3074
$ws.Cells.Item(79, 10).Value = 47621736
$ws.Cells.Item(79, 11).Value = 3074
$ws.Cells.Item(79, 12).Value = 47621736
$ws.Cells.Item(79, 13).Value = -1982
$ws.Cells.Item(79, 14).Value = -47623920

$ws.Cells.Item(111, 8).Value = 379.2381
$ws.Cells.Item(111, 9).Value = 260.30768
$ws.Cells.Item(111, 10).Value = 572.5
$ws.Cells.Item(111, 11).Value = 780.92304
$ws.Cells.Item(111, 12).Value = 1717.5
$ws.Cells.Item(111, 13).Value = 2286.07696
$ws.Cells.Item(111, 14).Value = -7851.5

$ws.Cells.Item(112, 8).Value = 989.2807
$ws.Cells.Item(112, 9).Value = 0
$ws.Cells.Item(112, 10).Value = 989.2807
$ws.Cells.Item(112, 11).Value = 0
$ws.Cells.Item(112, 12).Value = 2967.8421
$ws.Cells.Item(112, 13).ClearContents()
$ws.Cells.Item(112, 14).Value = -5183.8421

$ws.Cells.Item(115, 8).Value = 2088.2144
$ws.Cells.Item(115, 9).Value = 297
$ws.Cells.Item(115, 10).Value = 3083.3333
$ws.Cells.Item(115, 11).Value = 891
$ws.Cells.Item(115, 12).Value = 9249.999899999999
$ws.Cells.Item(115, 13).Value = 676

$ws.Cells.Item(116, 8).Value = 2507.4814
$ws.Cells.Item(116, 9).Value = 2501.3635
$ws.Cells.Item(116, 10).Value = 2511.6875
$ws.Cells.Item(116, 11).Value = 2501.3635
$ws.Cells.Item(116, 12).Value = 2511.6875
$ws.Cells.Item(116, 13).Value = 940.6365000000001
$ws.Cells.Item(116, 14).Value = -9395.6875

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(41, 8).Value = 1000
$ws.Cells.Item(41, 9).Value = 1000
$ws.Cells.Item(41, 10).Value = 0
$ws.Cells.Item(41, 11).Value = 1000
$ws.Cells.Item(41, 12).Value = 0
$ws.Cells.Item(41, 13).Value = -586

$ws.Cells.Item(45, 8).Value = 38182.332
$ws.Cells.Item(45, 9).Value = 59678.47
$ws.Cells.Item(45, 10).Value = 1638.9
$ws.Cells.Item(45, 11).Value = 59678.47
$ws.Cells.Item(45, 12).Value = 1638.9
$ws.Cells.Item(45, 13).Value = -59301.47
$ws.Cells.Item(45, 14).Value = -2392.9

$ws.Cells.Item(61, 8).Value = 4631234
$ws.Cells.Item(61, 9).Value = 5748703.5
$ws.Cells.Item(61, 10).Value = 1718.1428
$ws.Cells.Item(61, 11).Value = 5748703.5
$ws.Cells.Item(61, 12).Value = 1718.1428
$ws.Cells.Item(61, 13).Value = -5748491.5
$ws.Cells.Item(61, 14).Value = -2142.1428

$ws.Cells.Item(74, 8).Value = 1131.8246
$ws.Cells.Item(74, 9).Value = 1219.2162
$ws.Cells.Item(74, 10).Value = 970.15
$ws.Cells.Item(74, 11).Value = 1219.2162
$ws.Cells.Item(74, 12).Value = 970.15
$ws.Cells.Item(74, 13).Value = -345.2162000000001
$ws.Cells.Item(74, 14).Value = -2718.15

$ws.Cells.Item(77, 8).Value = 1131.8246
$ws.Cells.Item(77, 9).Value = 1219.2162
$ws.Cells.Item(77, 10).Value = 970.15
$ws.Cells.Item(77, 11).Value = 6096.081
$ws.Cells.Item(77, 12).Value = 4850.75
$ws.Cells.Item(77, 13).Value = -1728.081
$ws.Cells.Item(77, 14).Value = -13586.75

$ws.Cells.Item(88, 8).Value = 2382.1177
$ws.Cells.Item(88, 9).Value = 2242.5715
$ws.Cells.Item(88, 10).Value = 3033.3333
$ws.Cells.Item(88, 11).Value = 2242.5715
$ws.Cells.Item(88, 12).Value = 3033.3333
$ws.Cells.Item(88, 13).Value = -1836.5715
$ws.Cells.Item(88, 14).Value = -3845.3333

$ws.Cells.Item(91, 8).Value = 2382.1177
$ws.Cells.Item(91, 9).Value = 2242.5715
$ws.Cells.Item(91, 10).Value = 3033.3333
$ws.Cells.Item(91, 11).Value = 2242.5715
$ws.Cells.Item(91, 12).Value = 3033.3333
$ws.Cells.Item(91, 13).Value = -838.5715
$ws.Cells.Item(91, 14).Value = -5841.3333

$ws.Cells.Item(122, 8).Value = 1125.3334
$ws.Cells.Item(122, 9).Value = 1125.3334
$ws.Cells.Item(122, 10).Value = 0
$ws.Cells.Item(122, 11).Value = 3376.0002
$ws.Cells.Item(122, 12).Value = 0
$ws.Cells.Item(122, 13).Value = -926.0001999999999
$ws.Cells.Item(122, 14).ClearContents()

$ws.Cells.Item(136, 8).Value = 4631234
$ws.Cells.Item(136, 9).Value = 5748703.5
$ws.Cells.Item(136, 10).Value = 1718.1428
$ws.Cells.Item(136, 11).Value = 17246110.5
$ws.Cells.Item(136, 12).Value = 5154.428400000001
$ws.Cells.Item(136, 13).Value = -17243560.5
$ws.Cells.Item(136, 14).Value = -10254.4284

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 896702.75
$ws.Cells.Item(86, 9).Value = 2600.4614
$ws.Cells.Item(86, 10).Value = 1790805.1
$ws.Cells.Item(86, 11).Value = 2600.4614
$ws.Cells.Item(86, 12).Value = 1790805.1
$ws.Cells.Item(86, 13).Value = -1477.4614
$ws.Cells.Item(86, 14).Value = -1793051.1

$ws.Cells.Item(89, 8).Value = 896702.75
$ws.Cells.Item(89, 9).Value = 2600.4614
$ws.Cells.Item(89, 10).Value = 1790805.1
$ws.Cells.Item(89, 11).Value = 13002.307
$ws.Cells.Item(89, 12).Value = 8954025.5
$ws.Cells.Item(89, 13).Value = -7386.307000000001
$ws.Cells.Item(89, 14).Value = -8965257.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 1569
$ws.Cells.Item(31, 9).Value = 1176.0857
$ws.Cells.Item(31, 10).Value = 2078.3333
$ws.Cells.Item(31, 11).Value = 1176.0857
$ws.Cells.Item(31, 12).Value = 2078.3333
$ws.Cells.Item(31, 13).Value = -881.0857000000001
$ws.Cells.Item(31, 14).Value = -2668.3333

$ws.Cells.Item(34, 8).Value = 1569
$ws.Cells.Item(34, 9).Value = 1176.0857
$ws.Cells.Item(34, 10).Value = 2078.3333
$ws.Cells.Item(34, 11).Value = 1176.0857
$ws.Cells.Item(34, 12).Value = 2078.3333
$ws.Cells.Item(34, 13).Value = -974.0857000000001
$ws.Cells.Item(34, 14).Value = -2482.3333

$ws.Cells.Item(42, 8).Value = 45000
$ws.Cells.Item(42, 9).Value = 0
$ws.Cells.Item(42, 10).Value = 45000
$ws.Cells.Item(42, 11).Value = 0
$ws.Cells.Item(42, 12).Value = 45000
$ws.Cells.Item(42, 14).Value = -46186

$ws.Cells.Item(62, 8).Value = 3347.1428
$ws.Cells.Item(62, 9).Value = 2355.5715
$ws.Cells.Item(62, 10).Value = 4338.7144
$ws.Cells.Item(62, 11).Value = 2355.5715
$ws.Cells.Item(62, 12).Value = 4338.7144
$ws.Cells.Item(62, 13).Value = -1731.5715
$ws.Cells.Item(62, 14).Value = -5586.7144

$ws.Cells.Item(65, 8).Value = 3347.1428
$ws.Cells.Item(65, 9).Value = 2355.5715
$ws.Cells.Item(65, 10).Value = 4338.7144
$ws.Cells.Item(65, 11).Value = 11777.8575
$ws.Cells.Item(65, 12).Value = 21693.572
$ws.Cells.Item(65, 13).Value = -8657.8575
$ws.Cells.Item(65, 14).Value = -27933.572

$ws.Cells.Item(134, 8).Value = 953.55554
$ws.Cells.Item(134, 9).Value = 1057.3
$ws.Cells.Item(134, 10).Value = 657.1429000000001
$ws.Cells.Item(134, 11).Value = 3171.9
$ws.Cells.Item(134, 12).Value = 1971.4287
$ws.Cells.Item(134, 13).Value = -636.8999999999996
$ws.Cells.Item(134, 14).Value = -7041.4287

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(86, 8).Value = 200
$ws.Cells.Item(86, 9).Value = 200
$ws.Cells.Item(86, 10).Value = 0
$ws.Cells.Item(86, 11).Value = 600
$ws.Cells.Item(86, 12).Value = 0
$ws.Cells.Item(86, 13).Value = 586
$ws.Cells.Item(86, 14).ClearContents()

$ws.Cells.Item(89, 8).Value = 200
$ws.Cells.Item(89, 9).Value = 200
$ws.Cells.Item(89, 10).Value = 0
$ws.Cells.Item(89, 11).Value = 1800
$ws.Cells.Item(89, 12).Value = 0
$ws.Cells.Item(89, 13).Value = 4128
$ws.Cells.Item(89, 14).ClearContents()

$ws.Cells.Item(131, 8).Value = 901.6799999999999
$ws.Cells.Item(131, 9).Value = 499.5
$ws.Cells.Item(131, 10).Value = 909.88776
$ws.Cells.Item(131, 11).Value = 1498.5
$ws.Cells.Item(131, 12).Value = 2729.66328
$ws.Cells.Item(131, 13).Value = 3541.5
$ws.Cells.Item(131, 14).Value = -12809.66328

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(7, 8).Value = 3751249.8
$ws.Cells.Item(7, 9).Value = 5000000
$ws.Cells.Item(7, 10).Value = 2502499.5
$ws.Cells.Item(7, 11).Value = 5000000
$ws.Cells.Item(7, 12).Value = 2502499.5
$ws.Cells.Item(7, 13).Value = -4999888
$ws.Cells.Item(7, 14).Value = -2502723.5

$ws.Cells.Item(8, 8).Value = 3751249.8
$ws.Cells.Item(8, 9).Value = 5000000
$ws.Cells.Item(8, 10).Value = 2502499.5
$ws.Cells.Item(8, 11).Value = 5000000
$ws.Cells.Item(8, 12).Value = 2502499.5
$ws.Cells.Item(8, 13).Value = -4999861
$ws.Cells.Item(8, 14).Value = -2502777.5

$ws.Cells.Item(80, 8).Value = 7695673
$ws.Cells.Item(80, 9).Value = 3659.0908
$ws.Cells.Item(80, 10).Value = 50001750
$ws.Cells.Item(80, 11).Value = 3659.0908
$ws.Cells.Item(80, 12).Value = 50001750
$ws.Cells.Item(80, 13).Value = -2661.0908

$ws.Cells.Item(83, 8).Value = 7695673
$ws.Cells.Item(83, 9).Value = 3659.0908
$ws.Cells.Item(83, 10).Value = 50001750
$ws.Cells.Item(83, 11).Value = 18295.454
$ws.Cells.Item(83, 12).Value = 250008750
$ws.Cells.Item(83, 13).Value = -13303.454

$ws.Cells.Item(122, 8).Value = 15829.857
$ws.Cells.Item(122, 9).Value = 15829.857
$ws.Cells.Item(122, 10).Value = 0
$ws.Cells.Item(122, 11).Value = 47489.571
$ws.Cells.Item(122, 12).Value = 0
$ws.Cells.Item(122, 13).Value = -45039.571
$ws.Cells.Item(122, 14).ClearContents()

$ws.Cells.Item(132, 8).Value = 3616.739
$ws.Cells.Item(132, 9).Value = 1366.0278
$ws.Cells.Item(132, 10).Value = 11719.3
$ws.Cells.Item(132, 11).Value = 4098.0834
$ws.Cells.Item(132, 12).Value = 35157.89999999999
$ws.Cells.Item(132, 13).Value = -1568.0834
$ws.Cells.Item(132, 14).Value = -40217.89999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(50, 8).Value = 0
$ws.Cells.Item(50, 9).Value = 0
$ws.Cells.Item(50, 10).Value = 0
$ws.Cells.Item(50, 11).Value = 0
$ws.Cells.Item(50, 12).Value = 0
$ws.Cells.Item(50, 14).ClearContents()

$ws.Cells.Item(115, 8).Value = 29700
$ws.Cells.Item(115, 9).Value = 0
$ws.Cells.Item(115, 10).Value = 29700
$ws.Cells.Item(115, 11).Value = 0
$ws.Cells.Item(115, 12).Value = 29700
$ws.Cells.Item(115, 14).Value = -32050

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(2, 8).Value = 500000
$ws.Cells.Item(2, 9).Value = 500000
$ws.Cells.Item(2, 10).Value = 0
$ws.Cells.Item(2, 11).Value = 500000
$ws.Cells.Item(2, 12).Value = 0
$ws.Cells.Item(2, 13).Value = -499888

$ws.Cells.Item(81, 8).Value = 1334.3334
$ws.Cells.Item(81, 9).Value = 1220.6
$ws.Cells.Item(81, 10).Value = 1476.5
$ws.Cells.Item(81, 11).Value = 2441.2
$ws.Cells.Item(81, 12).Value = 2953
$ws.Cells.Item(81, 13).Value = -1380.2
$ws.Cells.Item(81, 14).Value = -5075

$ws.Cells.Item(84, 8).Value = 1334.3334
$ws.Cells.Item(84, 9).Value = 1220.6
$ws.Cells.Item(84, 10).Value = 1476.5
$ws.Cells.Item(84, 11).Value = 12206
$ws.Cells.Item(84, 12).Value = 14765
$ws.Cells.Item(84, 13).Value = -6902
$ws.Cells.Item(84, 14).Value = -25373
